# Generate Report for Handback
# Updates the localization-status workbook after a handback run:
#  - Status columns flip from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" populated for the handed-back file
#  - A hyperlink is added on the new "Latest Target File" cell, matching the one on
#    the "Source File Name" cell
#  - A few columns are widened to fit the newly-populated long filenames

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$sourceFileDisplay = "282bef07-6af4-4f2f-9dda-f4b4c87d8050.md"
$sourceFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/304477ec1bc0ed9375fb914ec6eb367706d93f76/e2e/282bef07-6af4-4f2f-9dda-f4b4c87d8050.md"

# ---- Overview sheet: zh-cn / de-de status columns ----
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# ---- zh-cn sheet: Status + handback info ----
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("I2").Value = $sourceFileDisplay
$zhcn.Range("J2").Value = "282bef07-6af4-4f2f-9dda-f4b4c87d8050.3ff8153dbe03c182fd803cb03642a53bbe5452d3.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-07 06:30:15"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $sourceFileUrl, "", "", $sourceFileDisplay)

# ---- de-de sheet: Status + handback info ----
$dede.Range("C2").Value = $statusText
$dede.Range("I2").Value = $sourceFileDisplay
$dede.Range("J2").Value = "282bef07-6af4-4f2f-9dda-f4b4c87d8050.3ff8153dbe03c182fd803cb03642a53bbe5452d3.de-de.xlf"
$dede.Range("K2").Value = "2016-09-07 06:30:41"
$dede.Hyperlinks.Add($dede.Range("I2"), $sourceFileUrl, "", "", $sourceFileDisplay)

# ---- Column widths: widen to fit the longer handback filenames ----
$overview.Columns.Item(5).ColumnWidth = 29.15   # E: zh-cn
$overview.Columns.Item(6).ColumnWidth = 29.15   # F: de-de

$zhcn.Columns.Item(3).ColumnWidth = 29.15    # C: Status
$zhcn.Columns.Item(9).ColumnWidth = 39.15    # I: Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.15   # J: Latest Handback File

$dede.Columns.Item(3).ColumnWidth = 29.15    # C: Status
$dede.Columns.Item(9).ColumnWidth = 39.15    # I: Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.15   # J: Latest Handback File
